# Renumber the "19.x" section titles to "20.x" across the deck.
# Each affected slide has a single title placeholder (ctrTitle) whose
# text run needs only its leading chapter number updated; all other
# formatting (font size, color, bold, etc.) stays untouched because
# setting TextRange.Text preserves the existing run formatting.

$p = $ppt.ActivePresentation

$titleUpdates = @{
    3  = "20.1 Mean-Shift Algorithm"
    5  = "20.2 Mean-Shift Example 1"
    6  = "20.2 Mean-Shift Example 1"
    7  = "20.2 Mean-Shift Example 1"
    8  = "20.2 Mean-Shift Example 2"
    9  = "20.2 Mean-Shift Example 2"
    10 = "20.2 Mean-Shift Example 2"
    11 = "20.3 Pros and Cons"
    12 = "20.3 Pros and Cons"
    13 = "20.3 Pros and Cons"
}

foreach ($slideIndex in $titleUpdates.Keys) {
    $s = $p.Slides.Item($slideIndex)
    $s.Shapes.Title.TextFrame.TextRange.Text = $titleUpdates[$slideIndex]
}
